$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array: A,B,C,D (strings) followed by E..T (numbers, as strings for exact precision)
$rowsData = @(
    @("FAPs", "Ly86", "Cd180", "ECs", "3", "1", "0.4933343333333333", "1.480003", "0.003704586072018129", "0.003704586072018129", "3", "1", "1.349047333333333", "4.047142", "0.0190104476758723", "0.0190104476758723", "0.665531366825111", "5.989782301426", "7.042583968286593E-05", "7.042583968286591E-05"),
    @("FAPs", "Ly86", "Cd180", "M2", "3", "1", "0.4933343333333333", "1.480003", "0.003704586072018129", "0.003704586072018129", "3", "1", "68.00131366666666", "204.003941", "0.9582580117159811", "0.958258011715981", "33.54738274353588", "301.926444691823", "0.003549949283602808", "0.003549949283602808"),
    @("FAPs", "Ly86", "Cd180", "sCs", "3", "1", "0.4933343333333333", "1.480003", "0.003704586072018129", "0.003704586072018129", "3", "1", "1.613109", "4.839327", "0.0227315406081467", "0.0227315406081467", "0.795802053109", "7.162218477981", "8.421094873245478E-05", "8.421094873245476E-05"),
    @("M2", "Ly86", "Cd180", "ECs", "3", "1", "130.9445773333333", "392.8337319999999", "0.9832996096535629", "0.983299609653563", "3", "1", "1.349047333333333", "4.047142", "0.0190104476758723", "0.0190104476758723", "176.6504328659937", "1589.853895793944", "0.01869296577902471", "0.01869296577902471"),
    @("M2", "Ly86", "Cd180", "M2", "3", "1", "130.9445773333333", "392.8337319999999", "0.9832996096535629", "0.983299609653563", "3", "1", "68.00131366666666", "204.003941", "0.9582580117159811", "0.958258011715981", "8904.403276193088", "80139.62948573779", "0.9422547288677235", "0.9422547288677235"),
    @("M2", "Ly86", "Cd180", "sCs", "3", "1", "130.9445773333333", "392.8337319999999", "0.9832996096535629", "0.983299609653563", "3", "1", "1.613109", "4.839327", "0.0227315406081467", "0.0227315406081467", "211.2278761975959", "1901.050885778364", "0.02235191500681476", "0.02235191500681476"),
    @("sCs", "Ly86", "Cd180", "ECs", "2", "0.6666666666666666", "1.730632333333333", "5.191897", "0.01299580427441884", "0.01299580427441884", "3", "1", "1.349047333333333", "4.047142", "0.0190104476758723", "0.0190104476758723", "2.334704934263778", "21.012344408374", "0.0002470560571647169", "0.0002470560571647169"),
    @("sCs", "Ly86", "Cd180", "M2", "2", "0.6666666666666666", "1.730632333333333", "5.191897", "0.01299580427441884", "0.01299580427441884", "3", "1", "68.00131366666666", "204.003941", "0.9582580117159811", "0.958258011715981", "117.6852721406752", "1059.167449266077", "0.01245333356465465", "0.01245333356465465"),
    @("sCs", "Ly86", "Cd180", "sCs", "2", "0.6666666666666666", "1.730632333333333", "5.191897", "0.01299580427441884", "0.01299580427441884", "3", "1", "1.613109", "4.839327", "0.0227315406081467", "0.0227315406081467", "2.791698592591", "25.125287333319", "0.0002954146525994783", "0.0002954146525994783")
)

for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $r = $i + 2
    $rowVals = $rowsData[$i]
    # Columns A-D: strings
    for ($c = 1; $c -le 4; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
    # Columns E-T: numbers (cast from exact string literal to preserve precision)
    for ($c = 5; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = [double]$rowVals[$c - 1]
    }
}